$d = $word.ActiveDocument
$r = $d.Paragraphs.Item(1).Range
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00365C75" w:rsidRPr="00814BB1" w:rsidRDefault="00911C96">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>Nome</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>=</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>${nome}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>Sobrenome = ${sobrenome}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>CPF = ${</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>cpf</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>Idade = ${idade}</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
          <w:color w:val="FF0000"/>
          <w:sz w:val="30"/>
          <w:szCs w:val="30"/>
        </w:rPr>
        <w:t>Nascimento = ${nascimento}</w:t>
      </w:r>

    </w:p>
'@
$r.InsertXML($newXml)
Write-Host "Done. ParaCount:" $d.Paragraphs.Count
